$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1983.8572
$ws.Range("J113").Value = 2231.3333
$ws.Range("L113").Value = 2231.3333
$ws.Range("N113").Value = -8739.3333
$ws.Range("H116").Value = 17429.285
$ws.Range("I116").Value = 19500.834
$ws.Range("K116").Value = 19500.834
$ws.Range("M116").Value = -16058.834
$ws.Range("H132").Value = 2709.3713
$ws.Range("I132").Value = 2488.4062
$ws.Range("K132").Value = 7465.2186
$ws.Range("M132").Value = -4935.2186
$ws.Range("H138").Value = 1878.5227
$ws.Range("I138").Value = 1360.2565
$ws.Range("J138").Value = 2291.0205
$ws.Range("K138").Value = 4080.7695
$ws.Range("L138").Value = 6873.0615
$ws.Range("M138").Value = 1059.2305
$ws.Range("N138").Value = -17153.0615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 36583.855
$ws.Range("I36").Value = 11500
$ws.Range("J36").Value = 70029
$ws.Range("K36").Value = 11500
$ws.Range("L36").Value = 70029
$ws.Range("M36").Value = -11154
$ws.Range("N36").Value = -70721
$ws.Range("H61").Value = 2367.152
$ws.Range("I61").Value = 1925.4138
$ws.Range("J61").Value = 3120.7058
$ws.Range("K61").Value = 1925.4138
$ws.Range("L61").Value = 3120.7058
$ws.Range("M61").Value = -1713.4138
$ws.Range("N61").Value = -3544.7058
$ws.Range("H104").Value = 74862.5
$ws.Range("J104").Value = 74862.5
$ws.Range("L104").Value = 74862.5
$ws.Range("N104").Value = -81850.5
$ws.Range("H106").Value = 41111
$ws.Range("J106").Value = 41111
$ws.Range("L106").Value = 41111
$ws.Range("N106").Value = -43635
$ws.Range("H132").Value = 2610.4866
$ws.Range("I132").Value = 2092.1897
$ws.Range("J132").Value = 4489.3125
$ws.Range("K132").Value = 6276.5691
$ws.Range("L132").Value = 13467.9375
$ws.Range("M132").Value = -3746.5691
$ws.Range("N132").Value = -18527.9375
$ws.Range("H135").Value = 22619.615
$ws.Range("J135").Value = 22619.615
$ws.Range("L135").Value = 22619.615
$ws.Range("N135").Value = -32759.615
$ws.Range("H136").Value = 2367.152
$ws.Range("I136").Value = 1925.4138
$ws.Range("J136").Value = 3120.7058
$ws.Range("K136").Value = 5776.2414
$ws.Range("L136").Value = 9362.117400000001
$ws.Range("M136").Value = -3226.2414
$ws.Range("N136").Value = -14462.1174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37094.11
$ws.Range("J81").Value = 37094.11
$ws.Range("L81").Value = 37094.11
$ws.Range("N81").Value = -39216.11
$ws.Range("H84").Value = 37094.11
$ws.Range("J84").Value = 37094.11
$ws.Range("L84").Value = 111282.33
$ws.Range("N84").Value = -121890.33
$ws.Range("H99").Value = 2269
$ws.Range("I99").Value = 2115
$ws.Range("K99").Value = 2115
$ws.Range("M99").Value = -617
$ws.Range("H134").Value = 3662.353
$ws.Range("I134").Value = 3085.3914
$ws.Range("J134").Value = 4868.727
$ws.Range("K134").Value = 9256.174199999999
$ws.Range("L134").Value = 14606.181
$ws.Range("M134").Value = -6721.174199999999
$ws.Range("N134").Value = -19676.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3254.4814
$ws.Range("I31").Value = 1049.1364
$ws.Range("J31").Value = 5877.054
$ws.Range("K31").Value = 1049.1364
$ws.Range("L31").Value = 5877.054
$ws.Range("M31").Value = -754.1364000000001
$ws.Range("N31").Value = -6467.054
$ws.Range("H34").Value = 3254.4814
$ws.Range("I34").Value = 1049.1364
$ws.Range("J34").Value = 5877.054
$ws.Range("K34").Value = 1049.1364
$ws.Range("L34").Value = 5877.054
$ws.Range("M34").Value = -847.1364000000001
$ws.Range("N34").Value = -6281.054
$ws.Range("H62").Value = 2824
$ws.Range("I62").Value = 2752.9033
$ws.Range("J62").Value = 3375
$ws.Range("K62").Value = 2752.9033
$ws.Range("L62").Value = 3375
$ws.Range("M62").Value = -2128.9033
$ws.Range("N62").Value = -4623
$ws.Range("H65").Value = 2824
$ws.Range("I65").Value = 2752.9033
$ws.Range("J65").Value = 3375
$ws.Range("K65").Value = 13764.5165
$ws.Range("L65").Value = 16875
$ws.Range("M65").Value = -10644.5165
$ws.Range("N65").Value = -23115
$ws.Range("H99").Value = 1941.1765
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H110").Value = 41111
$ws.Range("J110").Value = 41111
$ws.Range("L110").Value = 41111
$ws.Range("N110").Value = -49291
$ws.Range("H112").Value = 39999.332
$ws.Range("J112").Value = 39999.332
$ws.Range("L112").Value = 39999.332
$ws.Range("N112").Value = -42953.332
$ws.Range("H122").Value = 1995.2632
$ws.Range("I122").Value = 2222
$ws.Range("J122").Value = 1982.6666
$ws.Range("K122").Value = 6666
$ws.Range("L122").Value = 5947.9998
$ws.Range("M122").Value = -4216
$ws.Range("N122").Value = -10847.9998
$ws.Range("H123").Value = 39249.5
$ws.Range("J123").Value = 39249.5
$ws.Range("L123").Value = 39249.5
$ws.Range("N123").Value = -49049.5
$ws.Range("H126").Value = 1941.1765
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 1687.2778
$ws.Range("I132").Value = 1687.3103
$ws.Range("J132").Value = 1687.1428
$ws.Range("K132").Value = 5061.9309
$ws.Range("L132").Value = 5061.428400000001
$ws.Range("M132").Value = -2531.9309
$ws.Range("N132").Value = -10121.4284
$ws.Range("H134").Value = 4075.5264
$ws.Range("I134").Value = 4855.4287
$ws.Range("J134").Value = 1891.8
$ws.Range("K134").Value = 14566.2861
$ws.Range("L134").Value = 5675.4
$ws.Range("M134").Value = -12031.2861
$ws.Range("N134").Value = -10745.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 7095.7856
$ws.Range("J62").Value = 7095.7856
$ws.Range("L62").Value = 21287.3568
$ws.Range("N62").Value = -22659.3568
$ws.Range("H65").Value = 7095.7856
$ws.Range("J65").Value = 7095.7856
$ws.Range("L65").Value = 63862.0704
$ws.Range("N65").Value = -70726.0704
$ws.Range("H70").Value = 1337.3334
$ws.Range("I70").Value = 804.8
$ws.Range("K70").Value = 2414.4
$ws.Range("M70").Value = -2099.4
$ws.Range("H73").Value = 1337.3334
$ws.Range("I73").Value = 804.8
$ws.Range("K73").Value = 2414.4
$ws.Range("M73").Value = -1322.4
$ws.Range("H131").Value = 2851.4688
$ws.Range("J131").Value = 3240.018
$ws.Range("L131").Value = 9720.054
$ws.Range("N131").Value = -19800.054

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3949.9
$ws.Range("I132").Value = 3293.1765
$ws.Range("K132").Value = 9879.529500000001
$ws.Range("M132").Value = -7349.529500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2438.5264
$ws.Range("I132").Value = 2582.7878
$ws.Range("J132").Value = 2240.1667
$ws.Range("K132").Value = 7748.3634
$ws.Range("L132").Value = 6720.500100000001
$ws.Range("M132").Value = -5218.3634
$ws.Range("N132").Value = -11780.5001
$ws.Range("H136").Value = 5378234.5
$ws.Range("I136").Value = 2139.5
$ws.Range("J136").Value = 18519800
$ws.Range("K136").Value = 6418.5
$ws.Range("L136").Value = 55559400
$ws.Range("M136").Value = -3868.5
$ws.Range("N136").Value = -55564500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3050.8125
$ws.Range("I132").Value = 4145.8887
$ws.Range("J132").Value = 1642.8572
$ws.Range("K132").Value = 12437.6661
$ws.Range("L132").Value = 4928.571599999999
$ws.Range("M132").Value = -9907.666100000002
$ws.Range("N132").Value = -9988.571599999999
